# Updates the translator_testing_model workbook so the generator produces
# just 1 and 4 asset suite:
#  - TestAsset / AcceptanceTestAsset / TestEdgeData: move the
#    "test_runner_settings" column so it sits right after "test_metadata"
#    (immediately before "id"), shifting id/name/description/tags one
#    column to the right.
#  - TestCase / AcceptanceTestCase / QuantitativeTestCase: insert two new
#    columns ("components", "test_env") before the existing "id" column,
#    each with its own list data validation.

$wb = $excel.ActiveWorkbook

# --- Sheets that gain "components" + "test_env" columns just before "id" ---
$caseSheets = @("TestCase", "AcceptanceTestCase", "QuantitativeTestCase")

foreach ($name in $caseSheets) {
    $ws = $wb.Worksheets.Item($name)

    # "id" currently lives in column M; push it (and name/description/tags/
    # test_runner_settings after it) two columns to the right so the two
    # new columns can be written in place.
    $ws.Range("M1:N1").Insert(-4161) | Out-Null

    $ws.Range("M1").Value = "components"
    $ws.Range("N1").Value = "test_env"

    $componentsValidation = $ws.Range("M2:M1048576").Validation
    $componentsValidation.Add(3, 1, 1, '""') | Out-Null
    $componentsValidation.ShowInput = $false
    $componentsValidation.ShowError = $false

    $testEnvValidation = $ws.Range("N2:N1048576").Validation
    $testEnvValidation.Add(3, 1, 1, '"dev,ci,test,prod"') | Out-Null
    $testEnvValidation.ShowInput = $false
    $testEnvValidation.ShowError = $false
}

# --- Sheets where "test_runner_settings" moves next to "test_metadata" ---
$assetSheets = @("TestAsset", "AcceptanceTestAsset", "TestEdgeData")

foreach ($name in $assetSheets) {
    $ws = $wb.Worksheets.Item($name)

    # Locate "test_metadata" in row 1; "id","name","description","tags",
    # "test_runner_settings" occupy the five columns right after it.
    $metaCell = $ws.Rows(1).Find("test_metadata")
    $firstCol = $metaCell.Column + 1

    $idCol = $ws.Cells.Item(1, $firstCol)
    $nameCol = $ws.Cells.Item(1, $firstCol + 1)
    $descCol = $ws.Cells.Item(1, $firstCol + 2)
    $tagsCol = $ws.Cells.Item(1, $firstCol + 3)
    $settingsCol = $ws.Cells.Item(1, $firstCol + 4)

    $idCol.Value = "test_runner_settings"
    $nameCol.Value = "id"
    $descCol.Value = "name"
    $tagsCol.Value = "description"
    $settingsCol.Value = "tags"
}

Write-Output "done"
